$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: "NOM DU USE CASE : modifierquestionner" -> "NOM DU USE CASE : modifiersondage"
$ws.Range("A1").Value = "NOM DU USE CASE : modifiersondage"
$a1Bold = $ws.Range("A1").Characters(1, 15)
$a1Bold.Font.Bold = $true
$a1Bold.Font.Size = 12
$a1Rest = $ws.Range("A1").Characters(16, 18)
$a1Rest.Font.Bold = $false
$a1Rest.Font.Size = 11

# --- A2: "...BUT : modifier un questionnaire" -> "...BUT : modifier un sondage"
$ws.Range("A2").Value = "ACTEUR : Admin, super admin`n`nBUT : modifier un sondage"
$a2Acteur = $ws.Range("A2").Characters(1, 6)
$a2Acteur.Font.Bold = $true
$a2Acteur.Font.Size = 12
$a2AfterActeur = $ws.Range("A2").Characters(7, 23)
$a2AfterActeur.Font.Bold = $false
$a2AfterActeur.Font.Size = 11
$a2But = $ws.Range("A2").Characters(30, 3)
$a2But.Font.Bold = $true
$a2But.Font.Size = 12
$a2AfterBut = $ws.Range("A2").Characters(33, 22)
$a2AfterBut.Font.Bold = $false
$a2AfterBut.Font.Size = 11

# --- A6: "2) Quideance redirige vers la page de modification d'un questionnaire" -> "...d'un sondage"
$ws.Range("A6").Value = "2) Quideance redirige vers la page de modification d'un sondage"

# --- Update the active selection shown when the workbook is reopened
$null = $ws.Range("A6").Select()
